$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI values per row (only cells whose value actually changed)
$rowUpdates = @{
    2 = @{ "E" = 3; "F" = 1; "G" = 1.729584333333333; "H" = 5.188753; "I" = 0.2476387648475193; "J" = 0.2476387648475193; "M" = 50.86142466666666; "N" = 152.584274; "O" = 0.3434314568613803; "P" = 0.3434314568613804; "Q" = 87.96912327448022; "R" = 791.722109470322; "S" = 0.08504694178693632; "T" = 0.08504694178693634 }
    3 = @{ "E" = 3; "F" = 1; "G" = 1.729584333333333; "H" = 5.188753; "I" = 0.2476387648475193; "J" = 0.2476387648475193; "M" = 43.683024; "O" = 0.294960761928139; "P" = 0.294960761928139; "Q" = 75.553473943024; "R" = 679.981265487216; "S" = 0.07304371876236752; "T" = 0.07304371876236754 }
    4 = @{ "E" = 3; "F" = 1; "G" = 1.729584333333333; "H" = 5.188753; "I" = 0.2476387648475193; "J" = 0.2476387648475193; "M" = 36.64360566666667; "N" = 109.930817; "O" = 0.2474285170192034; "P" = 0.2474285170192035; "Q" = 63.37820627791123; "R" = 570.4038565012011; "S" = 0.06127289234268894; "T" = 0.06127289234268895 }
    5 = @{ "E" = 3; "F" = 1; "G" = 1.729584333333333; "H" = 5.188753; "I" = 0.2476387648475193; "J" = 0.2476387648475193; "M" = 16.90969166666667; "N" = 50.729075; "O" = 0.1141792641912772; "P" = 0.1141792641912772; "Q" = 29.24673778816389; "R" = 263.220640093475; "S" = 0.02827521195552648; "T" = 0.02827521195552648 }
    6 = @{ "I" = 0.2307941364328804; "J" = 0.2307941364328804; "M" = 50.86142466666666; "N" = 152.584274; "O" = 0.3434314568613803; "P" = 0.3434314568613804; "Q" = 81.98537838529623; "R" = 737.868405467666; "S" = 0.07926196651020828; "T" = 0.0792619665102083 }
    7 = @{ "I" = 0.2307941364328804; "J" = 0.2307941364328804; "M" = 43.683024; "O" = 0.294960761928139; "P" = 0.294960761928139; "Q" = 70.414253535472; "R" = 633.728281819248; "S" = 0.06807521433078925; "T" = 0.06807521433078927 }
    8 = @{ "I" = 0.2307941364328804; "J" = 0.2307941364328804; "M" = 36.64360566666667; "N" = 109.930817; "O" = 0.2474285170192034; "P" = 0.2474285170192035; "Q" = 59.06715935843923; "R" = 531.604434225953; "S" = 0.0571050509143153; "T" = 0.05710505091431531 }
    9 = @{ "I" = 0.2307941364328804; "J" = 0.2307941364328804; "M" = 16.90969166666667; "N" = 50.729075; "O" = 0.1141792641912772; "P" = 0.1141792641912772; "Q" = 27.25734638296389; "R" = 245.316117446675; "S" = 0.02635190467756753; "T" = 0.02635190467756753 }
    10 = @{ "G" = 2.743651333333334; "H" = 8.230954000000001; "I" = 0.3928310486309039; "J" = 0.3928310486309038; "M" = 50.86142466666666; "N" = 152.584274; "O" = 0.3434314568613803; "P" = 0.3434314568613804; "Q" = 139.5460156019329; "R" = 1255.914140417396; "S" = 0.1349105393316951; "T" = 0.1349105393316951 }
    11 = @{ "G" = 2.743651333333334; "H" = 8.230954000000001; "I" = 0.3928310486309039; "J" = 0.3928310486309038; "M" = 43.683024; "O" = 0.294960761928139; "P" = 0.294960761928139; "Q" = 119.850987041632; "R" = 1078.658883374688; "S" = 0.1158697454132012; "T" = 0.1158697454132012 }
    12 = @{ "G" = 2.743651333333334; "H" = 8.230954000000001; "I" = 0.3928310486309039; "J" = 0.3928310486309038; "M" = 36.64360566666667; "N" = 109.930817; "O" = 0.2474285170192034; "P" = 0.2474285170192035; "Q" = 100.5372775454909; "R" = 904.8354979094181; "S" = 0.09719760380184314; "T" = 0.09719760380184313 }
    13 = @{ "G" = 2.743651333333334; "H" = 8.230954000000001; "I" = 0.3928310486309039; "J" = 0.3928310486309038; "M" = 16.90969166666667; "N" = 50.729075; "O" = 0.1141792641912772; "P" = 0.1141792641912772; "Q" = 46.39429808750556; "R" = 417.54868278755; "S" = 0.04485316008416444; "T" = 0.04485316008416444 }
    14 = @{ "E" = 3; "F" = 1; "G" = 0.8991316666666668; "H" = 2.697395; "I" = 0.1287360500886965; "J" = 0.1287360500886965; "M" = 50.86142466666666; "N" = 152.584274; "O" = 0.3434314568613803; "P" = 0.3434314568613804; "Q" = 45.73111752958111; "R" = 411.58005776623; "S" = 0.04421200923254067; "T" = 0.04421200923254068 }
    15 = @{ "E" = 3; "F" = 1; "G" = 0.8991316666666668; "H" = 2.697395; "I" = 0.1287360500886965; "J" = 0.1287360500886965; "M" = 43.683024; "O" = 0.294960761928139; "P" = 0.294960761928139; "Q" = 39.27679017416; "R" = 353.49111156744; "S" = 0.03797208342178099; "T" = 0.037972083421781 }
    16 = @{ "E" = 3; "F" = 1; "G" = 0.8991316666666668; "H" = 2.697395; "I" = 0.1287360500886965; "J" = 0.1287360500886965; "M" = 36.64360566666667; "N" = 109.930817; "O" = 0.2474285170192034; "P" = 0.2474285170192035; "Q" = 32.94742623574611; "R" = 296.526836121715; "S" = 0.03185296996035607; "T" = 0.03185296996035607 }
    17 = @{ "E" = 3; "F" = 1; "G" = 0.8991316666666668; "H" = 2.697395; "I" = 0.1287360500886965; "J" = 0.1287360500886965; "M" = 16.90969166666667; "N" = 50.729075; "O" = 0.1141792641912772; "P" = 0.1141792641912772; "Q" = 15.20403925106945; "R" = 136.836353259625; "S" = 0.01469898747401877; "T" = 0.01469898747401878 }
}

foreach ($row in $rowUpdates.Keys) {
    foreach ($col in $rowUpdates[$row].Keys) {
        $ws.Range("$col$row").Value = $rowUpdates[$row][$col]
    }
}